$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'" + '29.528.78'
$ws.Range('E2').Value = '  +2.26%  '

$ws.Range('D3').Value = "'" + '1.986.87'
$ws.Range('E3').Value = '  +5.84%  '

$ws.Range('D4').Value = "'" + '1.000'
$ws.Range('E4').Value = '  -0.12%  '

$ws.Range('D5').Value = "'" + '328.72'
$ws.Range('E5').Value = '  +1.14%  '

$ws.Range('E6').Value = '  +0.18%  '

$ws.Range('D7').Value = "'" + '0.4685'
$ws.Range('E7').Value = '  +2.21%  '

$ws.Range('D8').Value = "'" + '0.3945'
$ws.Range('E8').Value = '  +1.87%  '

$ws.Range('D9').Value = "'" + '46.42'
$ws.Range('E9').Value = '  -0.20%  '

$ws.Range('D10').Value = "'" + '0.07966'
$ws.Range('E10').Value = '  +1.31%  '

$ws.Range('D11').Value = "'" + '1.003'
$ws.Range('E11').Value = '  +1.89%  '

$ws.Range('D12').Value = "'" + '22.91'
$ws.Range('E12').Value = '  +5.29%  '

$ws.Range('D13').Value = "'" + '1.980.06'
$ws.Range('E13').Value = '  +3.76%  '

$ws.Range('D14').Value = "'" + '7.263'
$ws.Range('E14').Value = '  +4.04%  '

$ws.Range('D15').Value = "'" + '5.881'
$ws.Range('E15').Value = '  +4.21%  '

$ws.Range('D16').Value = "'" + '0.07128'
$ws.Range('E16').Value = '  +2.38%  '

$ws.Range('D17').Value = "'" + '88.86'
$ws.Range('E17').Value = '  +0.95%  '

$ws.Range('D18').Value = "'" + '1.004'
$ws.Range('E18').Value = '  +0.16%  '

$ws.Range('D19').Value = "'" + '0.000009960'
$ws.Range('E19').Value = '  -0.04%  '

$ws.Range('D20').Value = "'" + '17.33'
$ws.Range('E20').Value = '  +2.26%  '

$ws.Range('D21').Value = "'" + '1.003'
$ws.Range('E21').Value = '  +0.04%  '

$ws.Range('D22').Value = "'" + '29.625.00'
$ws.Range('E22').Value = '  +2.53%  '

$ws.Range('D23').Value = "'" + '5.551'
$ws.Range('E23').Value = '  +5.77%  '

$ws.Range('D24').Value = "'" + '11.28'
$ws.Range('E24').Value = '  +2.95%  '

$ws.Range('B25').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C25').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D25').Value = "'" + '2.249.91'
$ws.Range('E25').Value = '  +5.72%  '

$ws.Range('B26').Value = 'Toncoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D26').Value = "'" + '2.119'
$ws.Range('E26').Value = '  +0.77%  '

$ws.Range('B27').Value = 'Monero'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D27').Value = "'" + '157.87'
$ws.Range('E27').Value = '  +1.15%  '

$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').Value = "'" + '19.63'
$ws.Range('E28').Value = '  +1.80%  '

$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').Value = "'" + '6.024'
$ws.Range('E29').Value = '  +0.23%  '

$ws.Range('B30').Value = 'BitcoinCash'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D30').Value = "'" + '120.35'
$ws.Range('E30').Value = '  +2.67%  '

$ws.Range('B31').Value = 'LidoDAOToken'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D31').Value = "'" + '1.963'
$ws.Range('E31').Value = '  +1.90%  '

$ws.Range('B32').Value = 'Stellar'
$ws.Range('C32').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D32').Value = "'" + '0.09457'
$ws.Range('E32').Value = '  +1.35%  '

$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').Value = "'" + '0.8922'
$ws.Range('E33').Value = '  -1.07%  '

$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').Value = "'" + '5.284'
$ws.Range('E34').Value = '  +0.60%  '

$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').Value = "'" + '1.348'
$ws.Range('E35').Value = '  +2.42%  '

$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').Value = "'" + '3.187'
$ws.Range('E36').Value = '  -2.11%  '

$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D37').Value = "'" + '0.05848'
$ws.Range('E37').Value = '  +1.67%  '

$ws.Range('B38').Value = 'TrustWalletToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D38').Value = "'" + '1.175'
$ws.Range('E38').Value = '  -0.61%  '

$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').Value = "'" + '0.02131'

$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').Value = "'" + '7.907'
$ws.Range('E40').Value = '  +3.09%  '

$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D41').Value = "'" + '0.5764'
$ws.Range('E41').Value = '  +2.11%  '

$ws.Range('D42').Value = "'" + '0.000003211'
$ws.Range('E42').Value = '  +100.62%  '

$ws.Range('B43').Value = 'Algorand'
$ws.Range('C43').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D43').Value = "'" + '0.1825'
$ws.Range('E43').Value = '  +3.54%  '

$ws.Range('B44').Value = 'Aptos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D44').Value = "'" + '9.832'
$ws.Range('E44').Value = '  +2.03%  '

$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = "'" + '12.15'
$ws.Range('E45').Value = '  +2.15%  '

$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').Value = "'" + '0.5384'
$ws.Range('E46').Value = '  +0.76%  '

$ws.Range('D47').Value = "'" + '2.677'
$ws.Range('E47').Value = '  +6.72%  '

$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').Value = "'" + '2.164'
$ws.Range('E48').Value = '  -4.04%  '

$ws.Range('D49').Value = "'" + '0.06965'
$ws.Range('E49').Value = '  -1.03%  '

$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D50').Value = "'" + '1.871'
$ws.Range('E50').Value = '  +1.60%  '

$ws.Range('B51').Value = 'Quant'
$ws.Range('C51').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D51').Value = "'" + '114.66'
$ws.Range('E51').Value = '  +1.55%  '
